# Update recalculated "Socio-economic capacity" (col I) and
# "Risk to well-being" (col J) result values for rows 2-36, per the
# refreshed model run referenced in the commit message
# ("rescues tot eq cost").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Cells.Item(2, 9).Value = 68.8088485414686
$ws.Cells.Item(2, 10).Value = 1.37636293020071
$ws.Cells.Item(3, 9).Value = 76.2393369099718
$ws.Cells.Item(3, 10).Value = 2.89042048217911
$ws.Cells.Item(4, 9).Value = 53.5058780289657
$ws.Cells.Item(4, 10).Value = 1.82286152615669
$ws.Cells.Item(5, 9).Value = 65.1196250282816
$ws.Cells.Item(5, 10).Value = 0.0430570084599718
$ws.Cells.Item(6, 9).Value = 130.239094175276
$ws.Cells.Item(6, 10).Value = 0.0114363588269259
$ws.Cells.Item(7, 9).Value = 59.9027816003588
$ws.Cells.Item(7, 10).Value = 0.254854614473707
$ws.Cells.Item(8, 9).Value = 174.394572287613
$ws.Cells.Item(8, 10).Value = 0.0513751402997716
$ws.Cells.Item(9, 9).Value = 95.7190425333085
$ws.Cells.Item(9, 10).Value = 0.715936799481044
$ws.Cells.Item(10, 9).Value = 64.0761237823329
$ws.Cells.Item(10, 10).Value = 0.737147723653262
$ws.Cells.Item(11, 9).Value = 90.7743402070683
$ws.Cells.Item(11, 10).Value = 0.711949475665127
$ws.Cells.Item(12, 9).Value = 68.6347293718318
$ws.Cells.Item(12, 10).Value = 0.238497180719832
$ws.Cells.Item(13, 9).Value = 74.5496670554051
$ws.Cells.Item(13, 10).Value = 0.45354473899707
$ws.Cells.Item(14, 9).Value = 127.226698069123
$ws.Cells.Item(14, 10).Value = 0.077831565311732
$ws.Cells.Item(15, 9).Value = 102.69382914028
$ws.Cells.Item(15, 10).Value = 0.0975713646128026
$ws.Cells.Item(16, 9).Value = 87.4199904204734
$ws.Cells.Item(16, 10).Value = 0.615357363498423
$ws.Cells.Item(17, 9).Value = 97.6439967885795
$ws.Cells.Item(17, 10).Value = 0.48645544297466
$ws.Cells.Item(18, 9).Value = 122.134674104081
$ws.Cells.Item(18, 10).Value = 0.0374181542255685
$ws.Cells.Item(19, 9).Value = 176.551295615766
$ws.Cells.Item(19, 10).Value = 0.066193926886424
$ws.Cells.Item(20, 9).Value = 66.6072730427029
$ws.Cells.Item(20, 10).Value = 0.116616626279725
$ws.Cells.Item(21, 9).Value = 38.569217916996
$ws.Cells.Item(21, 10).Value = 0.822724871710023
$ws.Cells.Item(22, 9).Value = 46.5946122751741
$ws.Cells.Item(22, 10).Value = 1.67457495501539
$ws.Cells.Item(23, 9).Value = 70.5199710452451
$ws.Cells.Item(23, 10).Value = 0.310252583596808
$ws.Cells.Item(24, 9).Value = 52.9553224441374
$ws.Cells.Item(24, 10).Value = 1.53451306486949
$ws.Cells.Item(25, 9).Value = 92.7204760333975
$ws.Cells.Item(25, 10).Value = 0.928616422288391
$ws.Cells.Item(26, 9).Value = 116.159677748139
$ws.Cells.Item(26, 10).Value = 0.679107427447446
$ws.Cells.Item(27, 9).Value = 159.443500953347
$ws.Cells.Item(27, 10).Value = 0.457178855250181
$ws.Cells.Item(28, 9).Value = 103.229184167326
$ws.Cells.Item(28, 10).Value = 0.490301827183283
$ws.Cells.Item(29, 9).Value = 97.1684707462559
$ws.Cells.Item(29, 10).Value = 0.0383885872937012
$ws.Cells.Item(30, 9).Value = 220.078240584589
$ws.Cells.Item(30, 10).Value = 0.0369039582945241
$ws.Cells.Item(31, 9).Value = 56.5861491156952
$ws.Cells.Item(31, 10).Value = 0.839062971768743
$ws.Cells.Item(32, 9).Value = 44.7773113589729
$ws.Cells.Item(32, 10).Value = 0.00342544615199732
$ws.Cells.Item(33, 9).Value = 88.3627851469908
$ws.Cells.Item(33, 10).Value = 0.127192757974572
$ws.Cells.Item(34, 9).Value = 47.4978727688213
$ws.Cells.Item(34, 10).Value = 1.66567230282632
$ws.Cells.Item(35, 9).Value = 112.958268860065
$ws.Cells.Item(35, 10).Value = 0.553447562294011
$ws.Cells.Item(36, 9).Value = 74.9023030223037
$ws.Cells.Item(36, 10).Value = 0.340884410709566
